$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("P2").Value = "'320018483205"
$ws.Range("P2").Style = "Normal"

# Row 3
$ws.Range("P3").Value = "'320018483238"
$ws.Range("P3").Style = "Normal"
$ws.Range("Q3").Value = "'$27.50"
$ws.Range("Q3").Style = "Normal"
$ws.Range("R3").Value = "PASS"

# Row 4
$ws.Range("P4").Value = "'320018483260"
$ws.Range("P4").Style = "Normal"
$ws.Range("Q4").Value = "'$31.73"
$ws.Range("Q4").Style = "Normal"
$ws.Range("R4").Value = "PASS"

# Row 5
$ws.Range("P5").Value = "'320018483282"
$ws.Range("P5").Style = "Normal"
$ws.Range("Q5").Value = "'$43.36"
$ws.Range("Q5").Style = "Normal"
$ws.Range("R5").Value = "PASS"

# Row 6
$ws.Range("P6").Value = "'320018483330"
$ws.Range("P6").Style = "Normal"
$ws.Range("Q6").Value = "'$56.05"
$ws.Range("Q6").Style = "Normal"
$ws.Range("R6").Value = "PASS"

# Row 7
$ws.Range("P7").Value = "'320018483352"
$ws.Range("P7").Style = "Normal"
$ws.Range("Q7").Value = "'$231.08"
$ws.Range("Q7").Style = "Normal"
$ws.Range("R7").Value = "PASS"

# Row 8
$ws.Range("P8").Value = "'320018483385"
$ws.Range("P8").Style = "Normal"
$ws.Range("Q8").Value = "'$19.04"
$ws.Range("Q8").Style = "Normal"
$ws.Range("R8").Value = "PASS"

# Row 9
$ws.Range("P9").Value = "'320018483411"
$ws.Range("P9").Style = "Normal"
$ws.Range("Q9").Value = "'$23.27"
$ws.Range("Q9").Style = "Normal"
$ws.Range("R9").Value = "PASS"

# Row 10
$ws.Range("P10").Value = "'320018483444"
$ws.Range("P10").Style = "Normal"
$ws.Range("Q10").Value = "'$27.50"
$ws.Range("Q10").Style = "Normal"
$ws.Range("R10").Value = "PASS"

# Row 11
$ws.Range("P11").Value = "'320018483466"
$ws.Range("P11").Style = "Normal"
$ws.Range("Q11").Value = "'$40.19"
$ws.Range("Q11").Style = "Normal"
$ws.Range("R11").Value = "PASS"

# Row 12
$ws.Range("P12").Value = "'320018483503"
$ws.Range("P12").Style = "Normal"
$ws.Range("Q12").Value = "'$52.88"
$ws.Range("Q12").Style = "Normal"
$ws.Range("R12").Value = "PASS"

# Row 13
$ws.Range("P13").Value = "'320018475505"
$ws.Range("P13").Style = "Normal"
$ws.Range("Q13").Value = "'$14.81"
$ws.Range("Q13").Style = "Normal"
$ws.Range("R13").Value = "PASS"

# Row 14
$ws.Range("P14").Value = "'320018475538"
$ws.Range("P14").Style = "Normal"
$ws.Range("Q14").Value = "'$17.98"
$ws.Range("Q14").Style = "Normal"
$ws.Range("R14").Value = "PASS"

# Row 15
$ws.Range("P15").Value = "'320018475550"
$ws.Range("P15").Style = "Normal"
$ws.Range("Q15").Value = "'$21.15"
$ws.Range("Q15").Style = "Normal"
$ws.Range("R15").Value = "PASS"

# Row 16
$ws.Range("P16").Value = "'320018475582"
$ws.Range("P16").Style = "Normal"
$ws.Range("Q16").Value = "'$31.73"
$ws.Range("Q16").Style = "Normal"
$ws.Range("R16").Value = "PASS"

# Row 17
$ws.Range("P17").Value = "'320018475696"
$ws.Range("P17").Style = "Normal"
$ws.Range("Q17").Value = "'$42.30"
$ws.Range("Q17").Style = "Normal"
$ws.Range("R17").Value = "PASS"

# Row 18
$ws.Range("P18").Value = "'320018475733"
$ws.Range("P18").Style = "Normal"
$ws.Range("Q18").Value = "'$43.36"
$ws.Range("Q18").Style = "Normal"
$ws.Range("R18").Value = "PASS"

# Row 19
$ws.Range("P19").Value = "'320018475766"
$ws.Range("P19").Style = "Normal"
$ws.Range("Q19").Value = "'$53.93"
$ws.Range("Q19").Style = "Normal"
$ws.Range("R19").Value = "PASS"

# Row 20
$ws.Range("P20").Value = "'320018475799"
$ws.Range("P20").Style = "Normal"
$ws.Range("Q20").Value = "'$62.39"
$ws.Range("Q20").Style = "Normal"
$ws.Range("R20").Value = "PASS"

# Row 21
$ws.Range("P21").Value = "'320018475836"
$ws.Range("P21").Style = "Normal"
$ws.Range("Q21").Value = "'$111.04"
$ws.Range("Q21").Style = "Normal"
$ws.Range("R21").Value = "PASS"

# Row 22
$ws.Range("P22").Value = "'320018475869"
$ws.Range("P22").Style = "Normal"
$ws.Range("Q22").Value = "'$223.37"
$ws.Range("Q22").Style = "Normal"
$ws.Range("R22").Value = "PASS"

# Row 23
$ws.Range("P23").Value = "'320018475870"
$ws.Range("P23").Style = "Normal"
$ws.Range("Q23").Value = "'$436.98"
$ws.Range("Q23").Style = "Normal"
$ws.Range("R23").Value = "PASS"

# Row 24
$ws.Range("P24").Value = "'320018475880"
$ws.Range("P24").Style = "Normal"
$ws.Range("Q24").Value = "'$278.12"
$ws.Range("Q24").Style = "Normal"
$ws.Range("R24").Value = "PASS"

# Row 25
$ws.Range("P25").Value = "'320018475891"
$ws.Range("P25").Style = "Normal"
$ws.Range("Q25").Value = "'$52.88"
$ws.Range("Q25").Style = "Normal"
$ws.Range("R25").Value = "PASS"

# Row 26
$ws.Range("P26").Value = "'320018475906"
$ws.Range("P26").Style = "Normal"
$ws.Range("Q26").Value = "'$1,171.41"
$ws.Range("Q26").Style = "Normal"
$ws.Range("R26").Value = "PASS"
$ws.Range("R26").Style = "Normal"
